$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$styleSave = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.989.00'
$ws.Range('D2').Style = $styleSave
$ws.Range('E2').Value = '  -0.50%  '
$styleSave = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.826.49'
$ws.Range('D3').Style = $styleSave
$ws.Range('E3').Value = '  +0.24%  '
$ws.Range('E4').Value = '  -0.39%  '
$styleSave = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.44'
$ws.Range('D5').Style = $styleSave
$ws.Range('E5').Value = '  +0.15%  '
$styleSave = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.006'
$ws.Range('D6').Style = $styleSave
$ws.Range('E6').Value = '  -0.33%  '
$styleSave = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4572'
$ws.Range('D7').Style = $styleSave
$ws.Range('E7').Value = '  -0.87%  '
$ws.Range('E8').Value = '  +1.87%  '
$ws.Range('E9').Value = '  +0.81%  '
$styleSave = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8747'
$ws.Range('D10').Style = $styleSave
$ws.Range('E10').Value = '  +0.52%  '
$styleSave = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07951'
$ws.Range('D11').Style = $styleSave
$ws.Range('E11').Value = '  +4.39%  '
$ws.Range('E12').Value = '  -1.67%  '
$styleSave = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.779.44'
$ws.Range('D13').Style = $styleSave
$ws.Range('E13').Value = '  -2.85%  '
$styleSave = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.587'
$ws.Range('D14').Style = $styleSave
$ws.Range('E14').Value = '  +1.61%  '
$styleSave = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.347'
$ws.Range('D15').Style = $styleSave
$ws.Range('E15').Value = '  +0.12%  '
$styleSave = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.52'
$ws.Range('D16').Style = $styleSave
$ws.Range('E16').Value = '  -0.93%  '
$ws.Range('E17').Value = '  -0.21%  '
$styleSave = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008929'
$ws.Range('D18').Style = $styleSave
$ws.Range('E18').Value = '  +3.41%  '
$styleSave = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.009'
$ws.Range('D19').Style = $styleSave
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$styleSave = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.82'
$ws.Range('D20').Style = $styleSave
$ws.Range('E20').Value = '  +2.39%  '
$ws.Range('B21').Value = 'WrappedBTC'
$ws.Range('C21').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$styleSave = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.448.62'
$ws.Range('D21').Style = $styleSave
$ws.Range('E21').Value = '  +0.05%  '
$styleSave = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.118'
$ws.Range('D22').Style = $styleSave
$ws.Range('E22').Value = '  -1.67%  '
$styleSave = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.173.73'
$ws.Range('D23').Style = $styleSave
$ws.Range('E23').Value = '  +3.63%  '
$ws.Range('E24').Value = '  -0.25%  '
$styleSave = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.31'
$ws.Range('D25').Style = $styleSave
$ws.Range('E25').Value = '  +1.20%  '
$styleSave = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.838'
$ws.Range('D26').Style = $styleSave
$ws.Range('E26').Value = '  -1.59%  '
$styleSave = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.39'
$ws.Range('D27').Style = $styleSave
$ws.Range('E27').Value = '  +1.02%  '
$styleSave = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.046'
$ws.Range('D28').Style = $styleSave
$ws.Range('E28').Value = '  -1.19%  '
$styleSave = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.151'
$ws.Range('D29').Style = $styleSave
$ws.Range('E29').Value = '  +0.98%  '
$styleSave = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.51'
$ws.Range('D30').Style = $styleSave
$ws.Range('E30').Value = '  -0.44%  '
$styleSave = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08878'
$ws.Range('D31').Style = $styleSave
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('E32').Value = '  +0.17%  '
$styleSave = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7297'
$ws.Range('D33').Style = $styleSave
$ws.Range('E33').Value = '  -0.47%  '
$styleSave = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.422'
$ws.Range('D34').Style = $styleSave
$ws.Range('E34').Value = '  -0.75%  '
$styleSave = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.130'
$ws.Range('D35').Style = $styleSave
$ws.Range('E35').Value = '  -0.46%  '
$styleSave = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.461'
$ws.Range('D36').Style = $styleSave
$ws.Range('E36').Value = '  -0.60%  '
$styleSave = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.071'
$ws.Range('D37').Style = $styleSave
$ws.Range('E37').Value = '  -0.17%  '
$styleSave = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01945'
$ws.Range('D38').Style = $styleSave
$ws.Range('E38').Value = '  +1.61%  '
$styleSave = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05227'
$ws.Range('D39').Style = $styleSave
$ws.Range('E39').Value = '  -0.34%  '
$styleSave = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.940'
$ws.Range('D40').Style = $styleSave
$ws.Range('E40').Value = '  +0.62%  '
$styleSave = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.149'
$ws.Range('D41').Style = $styleSave
$ws.Range('E41').Value = '  +0.11%  '
$styleSave = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5138'
$ws.Range('D42').Style = $styleSave
$ws.Range('E42').Value = '  -1.04%  '
$styleSave = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8858'
$ws.Range('D43').Style = $styleSave
$ws.Range('E43').Value = '  -12.33%  '
$styleSave = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1629'
$ws.Range('D44').Style = $styleSave
$ws.Range('E44').Value = '  +0.13%  '
$styleSave = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.182'
$ws.Range('D45').Style = $styleSave
$ws.Range('E45').Value = '  -1.16%  '
$styleSave = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4833'
$ws.Range('D46').Style = $styleSave
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('E47').Value = '  -0.31%  '
$styleSave = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.22'
$ws.Range('D48').Style = $styleSave
$ws.Range('E48').Value = '  +0.67%  '
$styleSave = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '102.70'
$ws.Range('D49').Style = $styleSave
$ws.Range('E49').Value = '  -0.73%  '
$styleSave = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.633'
$ws.Range('D50').Style = $styleSave
$ws.Range('E50').Value = '  -0.02%  '
$styleSave = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06199'
$ws.Range('D51').Style = $styleSave
$ws.Range('E51').Value = '  -0.99%  '
